# Auto-generated Excel COM-interop script
# Refreshes market-price-derived Leve profit columns (H-N) on several
# worksheets, matching a scheduled market-data re-pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 3767.3333
$ws.Range("J74").Value = 3999
$ws.Range("L74").Value = 3999
$ws.Range("N74").Value = -5871

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 3767.3333
$ws.Range("J77").Value = 3999
$ws.Range("L77").Value = 19995
$ws.Range("N77").Value = -29355

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2768.0625
$ws.Range("I132").Value = 2599.2903
$ws.Range("K132").Value = 7797.8709
$ws.Range("M132").Value = -5267.8709

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2600.8125
$ws.Range("I137").Value = 1203.8077
$ws.Range("J137").Value = 8654.5
$ws.Range("K137").Value = 3611.4231
$ws.Range("L137").Value = 25963.5
$ws.Range("M137").Value = -1061.4231
$ws.Range("N137").Value = -31063.5


$ws = $wb.Worksheets.Item("ARM")

# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1799.4706
$ws.Range("I2").Value = 1672.7333
$ws.Range("J2").Value = 2750
$ws.Range("K2").Value = 1672.7333
$ws.Range("L2").Value = 2750
$ws.Range("M2").Value = -1559.7333
$ws.Range("N2").Value = -2976

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1799.4706
$ws.Range("I116").Value = 1672.7333
$ws.Range("J116").Value = 2750
$ws.Range("K116").Value = 1672.7333
$ws.Range("L116").Value = 2750
$ws.Range("M116").Value = 621.2666999999999
$ws.Range("N116").Value = -7338


$ws = $wb.Worksheets.Item("BSM")

# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1799.4706
$ws.Range("I3").Value = 1672.7333
$ws.Range("J3").Value = 2750
$ws.Range("K3").Value = 1672.7333
$ws.Range("L3").Value = 2750
$ws.Range("M3").Value = -1558.7333
$ws.Range("N3").Value = -2978

# Row 75: I Saw the Pine / Hardsilver Saw
$ws.Range("H75").Value = 19774.5

# Row 78: I Came, I Sawed, I Conquered (L) / Hardsilver Saw
$ws.Range("H78").Value = 19774.5

# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 27597.8
$ws.Range("J82").Value = 43330
$ws.Range("L82").Value = 43330
$ws.Range("N82").Value = -44096

# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 27597.8
$ws.Range("J85").Value = 43330
$ws.Range("L85").Value = 43330
$ws.Range("N85").Value = -45982

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 2981.8572
$ws.Range("I99").Value = 2269.8
$ws.Range("J99").Value = 3377.4443
$ws.Range("K99").Value = 2269.8
$ws.Range("L99").Value = 3377.4443
$ws.Range("M99").Value = -771.8000000000002
$ws.Range("N99").Value = -6373.4443

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 32987.09
$ws.Range("J134").Value = 78262.21000000001
$ws.Range("L134").Value = 234786.63
$ws.Range("N134").Value = -239856.63


$ws = $wb.Worksheets.Item("CRP")

# Row 57: Clogs of War / Mahogany Pattens
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 2967.5715
$ws.Range("I99").Value = 3349
$ws.Range("J99").Value = 2014
$ws.Range("K99").Value = 3349
$ws.Range("L99").Value = 2014
$ws.Range("M99").Value = -1851
$ws.Range("N99").Value = -5010

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 2967.5715
$ws.Range("I126").Value = 3349
$ws.Range("J126").Value = 2014
$ws.Range("K126").Value = 10047
$ws.Range("L126").Value = 6042
$ws.Range("M126").Value = -7577
$ws.Range("N126").Value = -10982

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4557.8
$ws.Range("I132").Value = 4447.25
$ws.Range("K132").Value = 13341.75
$ws.Range("M132").Value = -10811.75

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1675159.4
$ws.Range("I134").Value = 3335500
$ws.Range("K134").Value = 10006500
$ws.Range("M134").Value = -10003965


$ws = $wb.Worksheets.Item("GSM")

# Row 53: North Ore South / Electrum Gorget
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 1577.4375
$ws.Range("I107").Value = 1039.1818
$ws.Range("J107").Value = 2761.6
$ws.Range("K107").Value = 1039.1818
$ws.Range("L107").Value = 2761.6
$ws.Range("M107").Value = 880.8181999999999
$ws.Range("N107").Value = -6601.6

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 45457920
$ws.Range("I132").Value = 47622450
$ws.Range("J132").Value = 2783
$ws.Range("K132").Value = 142867350
$ws.Range("L132").Value = 8349
$ws.Range("M132").Value = -142864820
$ws.Range("N132").Value = -13409


$ws = $wb.Worksheets.Item("LTW")

# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 11171749
$ws.Range("J7").Value = 174336.67
$ws.Range("L7").Value = 174336.67
$ws.Range("N7").Value = -174560.67

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 3131.8333
$ws.Range("I22").Value = 3599.5
$ws.Range("K22").Value = 3599.5
$ws.Range("M22").Value = -3304.5

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 3131.8333
$ws.Range("I27").Value = 3599.5
$ws.Range("K27").Value = 3599.5
$ws.Range("M27").Value = -3492.5

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3695.7666
$ws.Range("I46").Value = 3431.5334
$ws.Range("J46").Value = 3960
$ws.Range("K46").Value = 3431.5334
$ws.Range("L46").Value = 3960
$ws.Range("M46").Value = -3243.5334
$ws.Range("N46").Value = -4336

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 1499.5
$ws.Range("J68").Value = 999
$ws.Range("L68").Value = 999
$ws.Range("N68").Value = -2497

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 1499.5
$ws.Range("J71").Value = 999
$ws.Range("L71").Value = 4995
$ws.Range("N71").Value = -12483

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6117.1577
$ws.Range("I122").Value = 5185.5835
$ws.Range("J122").Value = 7714.143
$ws.Range("K122").Value = 15556.7505
$ws.Range("L122").Value = 23142.429
$ws.Range("M122").Value = -13106.7505
$ws.Range("N122").Value = -28042.429

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 11171749
$ws.Range("J126").Value = 174336.67
$ws.Range("L126").Value = 523010.01
$ws.Range("N126").Value = -527950.01

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 81632.92999999999
$ws.Range("I136").Value = 7080.5557
$ws.Range("J136").Value = 193461.5
$ws.Range("K136").Value = 21241.6671
$ws.Range("L136").Value = 580384.5
$ws.Range("M136").Value = -18691.6671
$ws.Range("N136").Value = -585484.5


$ws = $wb.Worksheets.Item("WVR")

# Row 5: Hire in the Blood / Hempen Halfgloves
$ws.Range("H5").Value = 5250037.5
$ws.Range("I5").Value = 500000
$ws.Range("J5").Value = 6833383.5
$ws.Range("K5").Value = 500000
$ws.Range("L5").Value = 6833383.5
$ws.Range("M5").Value = -499888
$ws.Range("N5").Value = -6833607.5

# Row 58: Seeing It Through to the End / Woolen Smock
$ws.Range("H58").Value = 32696.334
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 8344.079
$ws.Range("I132").Value = 1725.0416
$ws.Range("K132").Value = 5175.1248
$ws.Range("M132").Value = -2645.1248

